$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVA")

# Update row 3 (Q0)
$ws.Range("B3").Value = 0.1719895283108919
$ws.Range("C3").Value = 0.6769134546776546
$ws.Range("D3").Value = 0.7726047750576256
$ws.Range("E3").Value = 0.8789793939892024
$ws.Range("F3").Value = 0.8665135235198593
$ws.Range("G3").Value = 96

# Update row 4 (Q1)
$ws.Range("B4").Value = 0.09629569437071235
$ws.Range("C4").Value = 0.6864064840975582
$ws.Range("D4").Value = 0.6922636965734948
$ws.Range("E4").Value = 0.8320238557718732
$ws.Range("F4").Value = 0.8355647251700326
$ws.Range("G4").Value = 46
